$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 410; this shifts the existing rows
# 410..505 down to 411..506 (values, styles, and formatting all move
# with the rows), and the sheet's used range grows to R506.
$ws.Rows("410").Insert()

# Populate the newly inserted row 410 with the new weekly data point.
# Columns A, B, C, E, F, G, H, I and R carry the same constant metadata
# found on every data row in this sheet.
$ws.Range("A410").Value2 = 10
$ws.Range("B410").Value2 = "Vega Modelo de Temuco"
$ws.Range("C410").Value2 = "La Araucanía"
$ws.Range("D410").Value2 = 45173
$ws.Range("E410").Value2 = 9
$ws.Range("F410").Value2 = 100112001
$ws.Range("G410").Value2 = "Berenjena"
$ws.Range("H410").Value2 = "Sin especificar"
$ws.Range("I410").Value2 = "Primera"
$ws.Range("J410").Value2 = 100
$ws.Range("K410").Value2 = 10000
$ws.Range("L410").Value2 = 10000
$ws.Range("M410").Value2 = 10000
$ws.Range("N410").Value2 = "`$/caja 40 unidades"
$ws.Range("O410").Value2 = "Región de Arica y Parinacota"
$ws.Range("P410").Value2 = 250
$ws.Range("Q410").Value2 = 40
$ws.Range("R410").Value2 = "Hortaliza"

# Match the date-column number format used throughout column D.
$ws.Range("D410").NumberFormat = $ws.Range("D411").NumberFormat
